$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns (D: Price, E: Volume) keep their original text
# representation instead of being auto-coerced to numbers by Excels
# COM Value setter (e.g. "1.00" -> 1, "0.999" -> 0.999 as a float).
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.723.11"
$ws.Range("E2").Value = "  +2.65%  "
$ws.Range("D3").Value = "1.685.48"
$ws.Range("E3").Value = "  +2.93%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "220.43"
$ws.Range("E5").Value = "  +2.33%  "
$ws.Range("D6").Value = "0.528"
$ws.Range("E6").Value = "  +1.66%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").Value = "30.48"
$ws.Range("E8").Value = "  +5.50%  "
$ws.Range("D9").Value = "0.265"
$ws.Range("E9").Value = "  +2.17%  "
$ws.Range("D10").Value = "0.0627"
$ws.Range("E10").Value = "  +2.93%  "
$ws.Range("D11").Value = "0.0908"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").Value = "1.929.27"
$ws.Range("E12").Value = "  +2.98%  "
$ws.Range("D13").Value = "10.42"
$ws.Range("E13").Value = "  +12.17%  "
$ws.Range("D14").Value = "0.621"
$ws.Range("E14").Value = "  +9.29%  "
$ws.Range("D15").Value = "1.681.45"
$ws.Range("E15").Value = "  +2.84%  "
$ws.Range("E16").Value = "  +2.44%  "
$ws.Range("D17").Value = "30.730.21"
$ws.Range("E17").Value = "  +2.53%  "
$ws.Range("D18").Value = "66.37"
$ws.Range("E18").Value = "  +3.04%  "
$ws.Range("D19").Value = "246.24"
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("D20").Value = "0.0₃0716"
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "10.21"
$ws.Range("E22").Value = "  +2.38%  "
$ws.Range("E23").Value = "  +3.45%  "
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("D25").Value = "158.17"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  +1.67%  "
$ws.Range("D27").Value = "0.112"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("D28").Value = "6.70"
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").Value = "0.0499"
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("E32").Value = "  +2.61%  "
$ws.Range("D33").Value = "1.513.45"
$ws.Range("E33").Value = "  +5.75%  "
$ws.Range("E34").Value = "  +3.36%  "
$ws.Range("E35").Value = "  +4.90%  "
$ws.Range("D36").Value = "84.45"
$ws.Range("E36").Value = "  +9.83%  "
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("E38").Value = "  +4.15%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "0.588"
$ws.Range("E39").Value = "  +5.73%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.71"
$ws.Range("E40").Value = "  -5.02%  "
$ws.Range("D41").Value = "2.30"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("D43").Value = "1.99"
$ws.Range("E43").Value = "  -0.45%  "
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").Value = "52.10"
$ws.Range("E47").Value = "  -3.64%  "
$ws.Range("D48").Value = "1.820.50"
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("D49").Value = "5.44"
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("D50").Value = "94.77"
$ws.Range("E50").Value = "  +5.82%  "
$ws.Range("E51").Value = "  +1.45%  "
